# [PGE] Add 2021-01-22 data
# Update nombre_aides (C), nombre_entreprises (D), and montant_total (E)
# figures for several region rows to reflect the refreshed data extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 495872
$ws.Range("E2").Value = 764044962

$ws.Range("C19").Value = 117226
$ws.Range("E19").Value = 178790026

$ws.Range("C33").Value = 143111
$ws.Range("E33").Value = 222701540

$ws.Range("C45").Value = 103206
$ws.Range("E45").Value = 159030588

$ws.Range("C58").Value = 31011
$ws.Range("E58").Value = 50755569

$ws.Range("C67").Value = 216238
$ws.Range("E67").Value = 341218723

$ws.Range("C96").Value = 214593
$ws.Range("E96").Value = 323457894

$ws.Range("C111").Value = 857392
$ws.Range("E111").Value = 1399724697

$ws.Range("C152").Value = 132100
$ws.Range("D152").Value = 41347
$ws.Range("E152").Value = 206515908

$ws.Range("C164").Value = 350727
$ws.Range("E164").Value = 521029317

$ws.Range("C182").Value = 399651
$ws.Range("D182").Value = 122962
$ws.Range("E182").Value = 581422622

$ws.Range("E201").Value = 241519133

$ws.Range("C214").Value = 402275
$ws.Range("E214").Value = 605029400
